$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update pin assignments in column D (Pin_Out_MiniArcade sheet)
$ws.Range("D7").Value = "PA00"
$ws.Range("D8").Value = "PA01"
$ws.Range("D10").Value = "PA02"
$ws.Range("D27").Value = "PB09"
$ws.Range("D23").Value = "PB08"
$ws.Range("D19").Value = "PA03"
$ws.Range("D21").Value = "PB05"

# Update the view: scroll so A31 is top-left visible cell, and select I17
$ws.Range("I17").Select()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
